$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen the Description column (C) from 44 to 61 characters ---
$ws.Columns("C").ColumnWidth = 60 + 1/6

# --- Fill in Est (D) and By (E) for the three backlog item rows ---
$ws.Range("D4").Value = 24
$ws.Range("E4").Value = "Team 4"

$ws.Range("D5").Value = 24
$ws.Range("E5").Value = "Team 4"

$ws.Range("D6").Value = 24
$ws.Range("E6").Value = "Team 4"

# --- Row heights shrink now that column C is wider (less text wrapping) ---
$ws.Rows("4").RowHeight = 94.5
$ws.Rows("6").RowHeight = 90

# --- Clear the leftover Item # numbering in the now-empty rows 7 & 8 ---
$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()

# --- Remove the stray notes in column G ("User stories only" / "functional") ---
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# --- Trim the unused trailing blank rows 10-13 ---
$ws.Range("A10:E13").Delete()
